# Implement Meghan's data checks for non-breast solid malignancies
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 5
$ws.Range("L5").Value = 4.0
$ws.Range("M5").Value = 24.6

# Row 7
$ws.Range("L7").Value = 11.0
$ws.Range("M7").Value = 41.5

# Row 8
$ws.Range("L8").Value = 15.0
$ws.Range("M8").Value = 59.3

# Row 10
$ws.Range("L10").Value = 11.0
$ws.Range("M10").Value = 74.4

# Row 12
$ws.Range("L12").Value = 11.0
$ws.Range("M12").Value = 74.1

# Row 14 - fill previously empty cell F14
$ws.Range("F14").Value = 1800.0

# Row 23 - fill previously empty cells L23, M23
$ws.Range("L23").Value = 34.0
$ws.Range("M23").Value = 38.1

# Row 24 - fill previously empty cells L24, M24
$ws.Range("L24").Value = 31.0
$ws.Range("M24").Value = 34.7

# Row 81 - fill previously empty cells L81, M81
$ws.Range("L81").Value = 13.0
$ws.Range("M81").Value = 62.1

# Row 82 - fill previously empty cells L82, M82
$ws.Range("L82").Value = 4.0
$ws.Range("M82").Value = 19.3
